# Requirements for Rear View Camera.xlsx
# Adds new Hardware requirement rows (Description / Validation columns) for
# rows 7-13, applies center/wrap alignment formatting across the C:E
# requirement columns, widens columns D & E, and leaves the final selection
# on E19 (matching the author's saved state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# ---------------------------------------------------------------------
# 1. New requirement text (Description / Validation) for rows 7-13.
#    Written in the same order the shared-string table records them
#    (rows 7-10 in order, then the voltage-regulator validation text,
#    then rows 13/12/11 descriptions, then the row-11 validation text)
#    so new shared strings are interned in the original authoring order.
# ---------------------------------------------------------------------
$ws.Range("D7").Value = "The IR sensor will pass voltage ranges that an IC can read"
$ws.Range("E7").Value = "The IR Sensor will be tested using  a power supply and Multimeter.`nThe IR sensor will be connected for normal conditions, and move an object closer to the sensor to check the voltage change. Pass if the voltage change as the object moves closer and farther away is verified by the Multimeter"

$ws.Range("D8").Value = "The Receiver will be powered by two coin cells"
$ws.Range("E8").Value = "Receive will be able to turn on with two coin cells"

$ws.Range("D9").Value = "The Transmiter will be powered by two coin cells"
$ws.Range("E9").Value = "Transmitter will be able to turn on with two coin cells"

$ws.Range("D10").Value = "The speaker will be powered by two coin cells"
$ws.Range("E10").Value = "Speaker will be verified by spec sheets and test benching for loading "

$ws.Range("E12").Value = "with test bench equipment, passing a voltage less than the rated voltage, the voltage regulator will begin a smooth power down"
$ws.Range("D13").Value = "Transmitter will have voltage regulator sensing for smooth power down"
$ws.Range("D12").Value = "Receiver will have voltage regulator sensing for smooth power down"
$ws.Range("D11").Value = "IR sensor will have over current protection designed to maintain less than 1 A "
$ws.Range("E11").Value = "with test bench equipment, passing a 1.1 amp current will see the over current protection to activate"

$ws.Range("E13").Value = "with test bench equipment, passing a voltage less than the rated voltage, the voltage regulator will begin a smooth power down"

# ---------------------------------------------------------------------
# 2. Row heights for the newly-filled rows (taller to fit wrapped text)
# ---------------------------------------------------------------------
$ws.Rows.Item(7).RowHeight = 75
$ws.Rows.Item(11).RowHeight = 30
$ws.Rows.Item(12).RowHeight = 30
$ws.Rows.Item(13).RowHeight = 30

# ---------------------------------------------------------------------
# 3. Formatting: center horizontal/vertical alignment across the whole
#    requirements table. Column C is populated top-to-bottom (rows
#    6-92); columns D/E only hold data on the header rows (6/35/64) and
#    the new rows (7-13), so only those are touched - this mirrors the
#    sparse population the workbook actually has and avoids manufacturing
#    empty D/E cells on rows that never had any.
# ---------------------------------------------------------------------
$ws.Range("C6:C92").HorizontalAlignment = $xlCenter
$ws.Range("C6:C92").VerticalAlignment = $xlCenter

$deRange = $ws.Range("D6:E13")
$deRange.HorizontalAlignment = $xlCenter
$deRange.VerticalAlignment = $xlCenter

$ws.Range("D35:E35").HorizontalAlignment = $xlCenter
$ws.Range("D35:E35").VerticalAlignment = $xlCenter

$ws.Range("D64:E64").HorizontalAlignment = $xlCenter
$ws.Range("D64:E64").VerticalAlignment = $xlCenter

# Wrap text on the long free-form Description/Validation cells.
$ws.Range("D7:E7").WrapText = $true
$ws.Range("E11:E13").WrapText = $true

# ---------------------------------------------------------------------
# 4. Column widths (auto-widened by Excel to fit the new long text)
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 72.140625
$ws.Columns.Item(5).ColumnWidth = 70.85546875

# ---------------------------------------------------------------------
# 5. Leave the selection where the author last left it.
# ---------------------------------------------------------------------
$ws.Range("E19").Select() | Out-Null
